# edit.ps1 - apply quarterly financial update to WBS_QTR_FIN.xlsx
#
# Inserts two new quarterly columns (D, E) ahead of the existing quarters,
# shifting the previously reported quarters right by two columns (old D:K
# becomes F:M). The new quarters values are then populated and a handful
# of figures for the most recently reported quarter (now column F) are
# corrected to their revised amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E (existing D:K shifts right to F:M)
$ws.Range("D:E").EntireColumn.Insert()

# The newly inserted D:E columns have no explicit formatting yet. Copy the
# row-by-row formatting from column F (the former column D) into D:E, for
# just the row ranges that actually carry data/formatting, so every new
# cell matches the style used across its row (date format on the header
# rows, numeric format on data rows) without inventing cells on rows that
# only contain section-title text (5, 6, 37, 79).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the values for the two new quarters and correct the figures
# that were revised for the prior quarter (now in column F).
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 280400
$ws.Cells.Item(8, 5).Value = 268400
$ws.Cells.Item(8, 6).Value = 506400
$ws.Cells.Item(9, 4).Value = "NA"
$ws.Cells.Item(9, 5).Value = "NA"
$ws.Cells.Item(10, 4).Value = "NA"
$ws.Cells.Item(10, 5).Value = "NA"
$ws.Cells.Item(11, 4).Value = $null
$ws.Cells.Item(11, 5).Value = $null
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 4).Value = -1000
$ws.Cells.Item(15, 5).Value = -1000
$ws.Cells.Item(15, 6).Value = -1900
$ws.Cells.Item(16, 4).Value = $null
$ws.Cells.Item(16, 5).Value = $null
$ws.Cells.Item(17, 4).Value = 53300
$ws.Cells.Item(17, 5).Value = 48500
$ws.Cells.Item(17, 6).Value = 88700
$ws.Cells.Item(18, 4).Value = 227100
$ws.Cells.Item(18, 5).Value = 219900
$ws.Cells.Item(18, 6).Value = 417700
$ws.Cells.Item(19, 4).Value = $null
$ws.Cells.Item(19, 5).Value = $null
$ws.Cells.Item(20, 4).Value = -101600
$ws.Cells.Item(20, 5).Value = -106500
$ws.Cells.Item(20, 6).Value = -214900
$ws.Cells.Item(21, 4).Value = 135300
$ws.Cells.Item(21, 5).Value = 123200
$ws.Cells.Item(21, 6).Value = 221900
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 4).Value = 125500
$ws.Cells.Item(23, 5).Value = 113400
$ws.Cells.Item(23, 6).Value = 202700
$ws.Cells.Item(24, 4).Value = 37700
$ws.Cells.Item(24, 5).Value = 13700
$ws.Cells.Item(24, 6).Value = 40800
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 87900
$ws.Cells.Item(26, 5).Value = 99700
$ws.Cells.Item(26, 6).Value = 161900
$ws.Cells.Item(27, 4).Value = 85700
$ws.Cells.Item(27, 5).Value = 97500
$ws.Cells.Item(27, 6).Value = 157600
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 11000
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = 101600
$ws.Cells.Item(32, 5).Value = 106500
$ws.Cells.Item(32, 6).Value = 214900
$ws.Cells.Item(33, 4).Value = 96700
$ws.Cells.Item(33, 5).Value = 97500
$ws.Cells.Item(33, 6).Value = 157600
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 96700
$ws.Cells.Item(35, 5).Value = 97500
$ws.Cells.Item(35, 6).Value = 157600
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(39, 4).Value = $null
$ws.Cells.Item(39, 5).Value = $null
$ws.Cells.Item(40, 4).Value = $null
$ws.Cells.Item(40, 5).Value = $null
$ws.Cells.Item(41, 4).Value = 260400
$ws.Cells.Item(41, 5).Value = 222200
$ws.Cells.Item(42, 4).Value = 218400
$ws.Cells.Item(42, 5).Value = 233500
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 4).Value = 124900
$ws.Cells.Item(48, 5).Value = 128500
$ws.Cells.Item(49, 4).Value = 564100
$ws.Cells.Item(49, 5).Value = 565100
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 96500
$ws.Cells.Item(52, 5).Value = 92900
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 27610300
$ws.Cells.Item(54, 5).Value = 27346300
$ws.Cells.Item(55, 4).Value = $null
$ws.Cells.Item(55, 5).Value = $null
$ws.Cells.Item(56, 4).Value = $null
$ws.Cells.Item(56, 5).Value = $null
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(57, 5).Value = 0
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 0
$ws.Cells.Item(59, 4).Value = 230300
$ws.Cells.Item(59, 5).Value = 300200
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = 0
$ws.Cells.Item(61, 4).Value = 226000
$ws.Cells.Item(61, 5).Value = 226000
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(62, 5).Value = 0
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 24723800
$ws.Cells.Item(66, 5).Value = 24530100
$ws.Cells.Item(67, 4).Value = $null
$ws.Cells.Item(67, 5).Value = $null
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 145000
$ws.Cells.Item(70, 5).Value = 145000
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 1828300
$ws.Cells.Item(72, 5).Value = 1761000
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 2741500
$ws.Cells.Item(76, 5).Value = 2671200
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 96700
$ws.Cells.Item(81, 5).Value = 97500
$ws.Cells.Item(81, 6).Value = 157600
$ws.Cells.Item(82, 4).Value = $null
$ws.Cells.Item(82, 5).Value = $null
$ws.Cells.Item(83, 4).Value = 9800
$ws.Cells.Item(83, 5).Value = 9800
$ws.Cells.Item(83, 6).Value = 19200
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = -27000
$ws.Cells.Item(89, 5).Value = 239900
$ws.Cells.Item(89, 6).Value = 256500
$ws.Cells.Item(90, 4).Value = $null
$ws.Cells.Item(90, 5).Value = $null
$ws.Cells.Item(91, 4).Value = -7500
$ws.Cells.Item(91, 5).Value = -9400
$ws.Cells.Item(91, 6).Value = -16100
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -281300
$ws.Cells.Item(94, 5).Value = -426900
$ws.Cells.Item(94, 6).Value = -653200
$ws.Cells.Item(95, 4).Value = $null
$ws.Cells.Item(95, 5).Value = $null
$ws.Cells.Item(96, 4).Value = -30000
$ws.Cells.Item(96, 5).Value = -31000
$ws.Cells.Item(96, 6).Value = -54000
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = 346400
$ws.Cells.Item(100, 5).Value = 180700
$ws.Cells.Item(100, 6).Value = 394100
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(102, 4).Value = 38200
$ws.Cells.Item(102, 5).Value = -6400
$ws.Cells.Item(102, 6).Value = -2500
